$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.314.37"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "1.589.94"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  -0.46%  "
$ws.Range("D5").Value = "'211.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("E8").Value = "  +0.42%  "
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("E10").Value = "  -0.77%  "
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.595.36"
$ws.Range("E13").Value = "  +0.70%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'4.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("D16").Value = "'64.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").Value = "26.317.28"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("E19").Value = "  +3.73%  "
$ws.Range("D20").Value = "'212.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.43%  "
$ws.Range("E22").Value = "  +0.51%  "
$ws.Range("D23").Value = "'9.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.15%  "
$ws.Range("E24").Value = "  -2.94%  "
$ws.Range("D25").Value = "'144.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("E27").Value = "  +0.90%  "
$ws.Range("E28").Value = "  -0.93%  "
$ws.Range("D29").Value = "'15.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  +0.37%  "
$ws.Range("E32").Value = "  -0.72%  "
$ws.Range("D33").Value = "'2.99"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.54%  "
$ws.Range("D34").Value = "1.328.24"
$ws.Range("E34").Value = "  +3.45%  "
$ws.Range("E35").Value = "  -1.72%  "
$ws.Range("D36").Value = "'0.603"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.69%  "
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("E38").Value = "  -0.45%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("E40").Value = "  +4.91%  "
$ws.Range("E41").Value = "  -0.43%  "
$ws.Range("D42").Value = "'0.987"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -25.04%  "
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").Value = "'61.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.76%  "
$ws.Range("D46").Value = "1.725.86"
$ws.Range("E46").Value = "  +0.34%  "
$ws.Range("D47").Value = "'88.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.77%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.65%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.0504"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.10%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.0978"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.34%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "'1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.41%  "
